# Apply model-reduction database rebase fixes to the db_kpi export sheet.
# - Refreshes the "db timeofrun" timestamp and several computed KPI columns
#   for the existing runs (rows 2-7).
# - Appends the results of a second simulation run (rows 8-13) that was
#   missing because of the non-relative db_carrides path bug.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual KPI cells on the already-exported rows (2-7) ---

# Row 2 (ebuurt id 0 - rijtjeshuizen (laag))
$ws.Cells.Item(2, 2).Value  = 44715.54143104167     # B2  db timeofrun
$ws.Cells.Item(2, 20).Value = 16.592177558161886    # T2  averagetrafoload pct
$ws.Cells.Item(2, 21).Value = 0.739408956016037     # U2  maxtrafoloadin fr
$ws.Cells.Item(2, 22).Value = 0.6396518258887574    # V2  maxtrafoloadout fr
$ws.Cells.Item(2, 23).Value = 6800.0                # W2  totalcostmonthly eurpm
$ws.Cells.Item(2, 24).Value = 317.7                 # X2  co2emission ton
$ws.Cells.Item(2, 25).Value = 55.15491297789812     # Y2  selfconsumption pct

# Row 3 (ebuurt id 1 - VVD-wijk (laag))
$ws.Cells.Item(3, 2).Value  = 44715.541431261576    # B3
$ws.Cells.Item(3, 20).Value = 17.36224718755243     # T3
$ws.Cells.Item(3, 21).Value = 1.2207691706471688    # U3
$ws.Cells.Item(3, 22).Value = 0.26595531990458554   # V3
$ws.Cells.Item(3, 23).Value = 8200.0                # W3
$ws.Cells.Item(3, 24).Value = 296.77                # X3
$ws.Cells.Item(3, 25).Value = 86.17726211857631     # Y3

# Row 4 (ebuurt id 2 - rijtjeshuizen (hoog))
$ws.Cells.Item(4, 2).Value  = 44715.54143131944     # B4
$ws.Cells.Item(4, 20).Value = 25.622597085337024    # T4
$ws.Cells.Item(4, 21).Value = 1.818451052098019     # U4
$ws.Cells.Item(4, 23).Value = 5800.0                # W4
$ws.Cells.Item(4, 24).Value = 198.71                # X4

# Row 5 (ebuurt id 3 - VVD-wijk (hoog))
$ws.Cells.Item(5, 2).Value  = 44715.54143143519     # B5
$ws.Cells.Item(5, 20).Value = 35.72940059895572     # T5
$ws.Cells.Item(5, 21).Value = 3.272610974059578     # U5
$ws.Cells.Item(5, 23).Value = 7200.0                # W5
$ws.Cells.Item(5, 24).Value = 111.04                # X5

# Row 6 (holon 0 - WindHolon)
$ws.Cells.Item(6, 2).Value  = 44715.54143145833     # B6
$ws.Cells.Item(6, 34).Value = 155.65267225175722    # AH6 householdaveragemonthlyopex eurphpm
$ws.Cells.Item(6, 35).Value = 47.16187409387012     # AI6 selfconsumption pct
$ws.Cells.Item(6, 36).Value = 58.15151520670435     # AJ6 selfsufficiency pct

# Row 7 (main/root totals)
$ws.Cells.Item(7, 2).Value  = 44715.54143152778     # B7
$ws.Cells.Item(7, 38).Value = 46.211079155152746    # AL7 nationalco2households mton
$ws.Cells.Item(7, 39).Value = 47743.0               # AM7 nationalelhouseholdelectricityconsumption gwh

# Apply the same date format used by the other "db timeofrun" cells
$ws.Range("B2:B7").NumberFormat = "yyyy-mm-dd"

# --- Append the second simulation run (rows 8-13) ---

$newRows = @{
    8  = @(2.0, 44715.54306967593, 0.0, "pop_buurten[0]", "rijtjeshuizen (laag)", 40.0, 40.0, 0.0, 15.0, 0.0, 27.5, 55.0, 2.5, 0.0, 0.1, 95.0, 0.0, 5.0, 0.0, 16.045884873030943, 0.6845742062270038, 0.6396518258887574, 6600.0, 330.6, 51.92390300909424, -2.5859575815597293, 4.05959872840634, -5.858063759612174, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)

    9  = @(2.0, 44715.54306986111, 1.0, "pop_buurten[1]", "VVD-wijk (laag)", 40.0, 22.5, 0.0, 62.5, 22.5, 0.0, 0.0, 15.0, 0.0, 0.1, 92.5, 0.0, 7.5, 0.0, 14.14452876034304, 0.6191124122990689, 0.2739653785483065, 7700.0, 341.2, 77.4444014570355, -5.9241459786379025, 14.969228292534202, -10.13360188854081, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)

    10  = @(2.0, 44715.54306991898, 2.0, "pop_buurten[2]", "rijtjeshuizen (hoog)", 40.0, 87.5, 0.0, 35.0, 0.0, 12.5, 52.5, 0.0, 0.0, 0.225, 0.0, 0.0, 12.5, 87.5, 21.038341721770006, 0.8699073342047108, -0.0, 5300.0, 270.74, 100.0, -8.322974863051812, 36.24300935079269, 0.0, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)

    11  = @(2.0, 44715.543069976855, 3.0, "pop_buurten[3]", "VVD-wijk (hoog)", 40.0, 77.5, 0.0, 72.5, 22.5, 0.0, 0.0, 5.0, 0.0, 0.15, 0.0, 0.0, 15.0, 85.0, 24.82182455346574, 1.5075240486917962, -0.0, 5900.0, 263.56, 100.0, -17.423517264573913, 137.35122953144509, 0.0, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)

    12  = @(2.0, 44715.543070011576, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 0.0, "holonAgent[0]", "WindHolon", 70.0, 69.0, 137.83915743061382, 30.55314281402029, 51.20810434587503, $null, $null, $null)

    13  = @(2.0, 44715.54307006944, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, "root", 60.30465406439814, 40439.0)
}

foreach ($r in 8..13) {
    $rowValues = $newRows[$r]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $val = $rowValues[$col - 1]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
    $ws.Cells.Item($r, 2).NumberFormat = "yyyy-mm-dd"
}
